# This workbook contains a daily/weekly log of produce price observations
# (one row per market day). A new weekly observation is inserted at row 54,
# pushing all subsequent rows down by one (the former row 172 becomes row 173).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 54, shifting rows 54:172 down to 55:173.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new observation.
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44868
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112001
$ws.Range("G54").Value = "Berenjena"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 500
$ws.Range("K54").Value = 11500
$ws.Range("L54").Value = 12000
$ws.Range("M54").Value = 11750
$ws.Range("N54").Value = "`$/caja 40 unidades"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 294
$ws.Range("Q54").Value = 40
$ws.Range("R54").Value = "Hortaliza"
